# Fix Training Data Issue (#48)
#
# The "Date" column on the sheet held the literal string "4-25-2007-08"
# (a mash-up of the source file's "4-25-2007-08" name) for every data row.
# The box scores in this sheet are actually from 2008-04-25 - NBA stats
# for the 2007-08 season were captured one calendar day off from the real
# game date. Correct every data row to the real ISO date string
# "2008-04-25", keeping it stored as literal text (not letting Excel
# reinterpret the assignment as a date literal / numeric serial).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "4-25-2007-08"
$newValue = "2008-04-25"

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1
$firstCol = $used.Column
$lastCol = $firstCol + $used.Columns.Count - 1

# Locate the "Date" column header on row 1 (BF in the current layout), so
# the fix keeps working even if the sheet layout shifts.
$dateCol = 0
for ($col = $firstCol; $col -le $lastCol; $col++) {
    if ($ws.Cells.Item(1, $col).Value2 -eq "Date") {
        $dateCol = $col
        break
    }
}
if ($dateCol -eq 0) {
    $dateCol = 58  # fallback: column BF
}

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $dateCol)
    if ($cell.Value2 -eq $oldValue) {
        # A bare assignment of "2008-04-25" would be auto-parsed by Excel
        # as a date literal and stored as a numeric serial instead of text.
        # Prefix with an apostrophe to force literal-text entry (matching
        # the source file, which stores this column as plain text), then
        # restore the cell's original (unstyled) appearance.
        $cell.Value = "'" + $newValue
        $cell.Style = "Normal"
    }
}
